$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.692.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.89%  "

$ws.Range("D3").Value = "'2.617.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.36%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'328.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("D6").Value = "'110.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.24%  "

$ws.Range("D7").Value = "'0.534"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "'0.558"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.26%  "

$ws.Range("D10").Value = "'40.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.20%  "

$ws.Range("D11").Value = "'20.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.17%  "

$ws.Range("D12").Value = "'0.0820"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").Value = "'7.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").Value = "'3.030.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.55%  "

$ws.Range("D16").Value = "'2.620.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.41%  "

$ws.Range("D17").Value = "'0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.32%  "

$ws.Range("D18").Value = "'49.671.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "

$ws.Range("D19").Value = "'3.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.39%  "

$ws.Range("D20").Value = "'13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").Value = "'6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").Value = "'0.0₃0952"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "'280.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'72.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").Value = "'2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").Value = "'26.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.21%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").Value = "'36.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "

$ws.Range("D32").Value = "'49.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("D33").Value = "'19.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("D34").Value = "'5.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").Value = "'0.0793"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("E37").Value = "  +5.12%  "

$ws.Range("E38").Value = "  +1.72%  "

$ws.Range("E39").Value = "  +5.72%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.112"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'123.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'22.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("D44").Value = "'0.0315"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.19%  "

$ws.Range("D45").Value = "'3.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.65%  "

$ws.Range("D46").Value = "'2.053.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  +11.83%  "

$ws.Range("E48").Value = "  +9.22%  "

$ws.Range("D49").Value = "'9.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("D50").Value = "'5.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.42%  "

$ws.Range("D51").Value = "'81.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.51%  "
